$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(354).Insert()
$ws.Rows.Item(355).Insert()

$ws.Range("A354").Value = 1
$ws.Range("B354").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C354").Value = "Arica y Parinacota"
$ws.Range("D354").Value = 44918
$ws.Range("E354").Value = 15
$ws.Range("F354").Value = 100112043
$ws.Range("G354").Value = "Pepino ensalada"
$ws.Range("H354").Value = "Sin especificar"
$ws.Range("I354").Value = "Primera"
$ws.Range("J354").Value = 450
$ws.Range("K354").Value = 13000
$ws.Range("L354").Value = 14000
$ws.Range("M354").Value = 13556
$ws.Range("N354").Value = "`$/caja 70 unidades"
$ws.Range("O354").Value = "Región de Arica y Parinacota"
$ws.Range("P354").Value = 194
$ws.Range("Q354").Value = 70
$ws.Range("R354").Value = "Hortaliza"

$ws.Range("A355").Value = 1
$ws.Range("B355").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C355").Value = "Arica y Parinacota"
$ws.Range("D355").Value = 44918
$ws.Range("E355").Value = 15
$ws.Range("F355").Value = 100112043
$ws.Range("G355").Value = "Pepino ensalada"
$ws.Range("H355").Value = "Sin especificar"
$ws.Range("I355").Value = "Segunda"
$ws.Range("J355").Value = 190
$ws.Range("K355").Value = 11000
$ws.Range("L355").Value = 12000
$ws.Range("M355").Value = 11474
$ws.Range("N355").Value = "`$/caja 100 unidades"
$ws.Range("O355").Value = "Región de Arica y Parinacota"
$ws.Range("P355").Value = 115
$ws.Range("Q355").Value = 100
$ws.Range("R355").Value = "Hortaliza"
